# Apply updated loading_percent values (Case 380 kV) to B2:O25
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,14
$arr[0,0] = 11.01197132713819
$arr[0,1] = 7.67171550768592
$arr[0,2] = 0
$arr[0,3] = 21.78552987994918
$arr[0,4] = 37.87772308245467
$arr[0,5] = 3.609602093350006
$arr[0,6] = 0
$arr[0,7] = 18.19053903258717
$arr[0,8] = 7.57147938089041
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 18.87438220162424
$arr[0,12] = 16.87013181718938
$arr[0,13] = 19.29644185346895
$arr[1,0] = 10.48645652526517
$arr[1,1] = 7.262507498494661
$arr[1,2] = 0
$arr[1,3] = 21.75493604695122
$arr[1,4] = 37.80078641722215
$arr[1,5] = 3.61149475421795
$arr[1,6] = 0
$arr[1,7] = 18.28675726608654
$arr[1,8] = 7.587567856043336
$arr[1,9] = 0
$arr[1,10] = 0
$arr[1,11] = 18.68933900617738
$arr[1,12] = 16.91646567094952
$arr[1,13] = 19.3445445993421
$arr[2,0] = 10.15071500087388
$arr[2,1] = 6.997673533625544
$arr[2,2] = 0
$arr[2,3] = 21.74024970239971
$arr[2,4] = 37.76380188521851
$arr[2,5] = 3.612718507496169
$arr[2,6] = 0
$arr[2,7] = 18.35008256406965
$arr[2,8] = 7.597941711227123
$arr[2,9] = 0
$arr[2,10] = 0
$arr[2,11] = 18.57792433885727
$arr[2,12] = 16.94666956045301
$arr[2,13] = 19.37941940028432
$arr[3,0] = 10.01078040189119
$arr[3,1] = 6.886392702329099
$arr[3,2] = 0
$arr[3,3] = 21.73529955997374
$arr[3,4] = 37.75131865847059
$arr[3,5] = 3.613232748822579
$arr[3,6] = 0
$arr[3,7] = 18.37695419673884
$arr[3,8] = 7.602294102152234
$arr[3,9] = 0
$arr[3,10] = 0
$arr[3,11] = 18.5331196019674
$arr[3,12] = 16.95941994663899
$arr[3,13] = 19.39496897225789
$arr[4,0] = 9.98736159972618
$arr[4,1] = 6.867713602942603
$arr[4,2] = 0
$arr[4,3] = 21.73454018883448
$arr[4,4] = 37.74940238081747
$arr[4,5] = 3.613319078944476
$arr[4,6] = 0
$arr[4,7] = 18.38148054485042
$arr[4,8] = 7.603024371598285
$arr[4,9] = 0
$arr[4,10] = 0
$arr[4,11] = 18.52571716357702
$arr[4,12] = 16.96156386200017
$arr[4,13] = 19.39763161758345
$arr[5,0] = 10.14884017082557
$arr[5,1] = 6.996186278973644
$arr[5,2] = 0
$arr[5,3] = 21.74017874889876
$arr[5,4] = 37.76362304162537
$arr[5,5] = 3.61272537969958
$arr[5,6] = 0
$arr[5,7] = 18.3504406503106
$arr[5,8] = 7.597999902640151
$arr[5,9] = 0
$arr[5,10] = 0
$arr[5,11] = 18.57731760965934
$arr[5,12] = 16.94683972566279
$arr[5,13] = 19.37962369724967
$arr[6,0] = 10.83358125027122
$arr[6,1] = 7.533468613026765
$arr[6,2] = 0
$arr[6,3] = 21.77413316772583
$arr[6,4] = 37.84907347673193
$arr[6,5] = 3.610241916563903
$arr[6,6] = 0
$arr[6,7] = 18.22283250365508
$arr[6,8] = 7.576924077598571
$arr[6,9] = 0
$arr[6,10] = 0
$arr[6,11] = 18.81014873960444
$arr[6,12] = 16.88574402826277
$arr[6,13] = 19.31191665520469
$arr[7,0] = 12.06641085832628
$arr[7,1] = 8.477446891040588
$arr[7,2] = 0
$arr[7,3] = 21.87302791002638
$arr[7,4] = 38.09746784698005
$arr[7,5] = 3.605858792161246
$arr[7,6] = 0
$arr[7,7] = 18.00637493806664
$arr[7,8] = 7.539508627177593
$arr[7,9] = 0
$arr[7,10] = 0
$arr[7,11] = 19.28223851052151
$arr[7,12] = 16.7798218138368
$arr[7,13] = 19.22171089677124
$arr[8,0] = 12.89834844998137
$arr[8,1] = 9.102253936349737
$arr[8,2] = 0
$arr[8,3] = 21.96506122409031
$arr[8,4] = 38.32833864081745
$arr[8,5] = 3.60293222356727
$arr[8,6] = 0
$arr[8,7] = 17.86805482928834
$arr[8,8] = 7.514381476507534
$arr[8,9] = 0
$arr[8,10] = 0
$arr[8,11] = 19.63585696960525
$arr[8,12] = 16.71041597980399
$arr[8,13] = 19.18163672589195
$arr[9,0] = 13.25972238978505
$arr[9,1] = 9.371311849801241
$arr[9,2] = 0
$arr[9,3] = 22.01105320303305
$arr[9,4] = 38.44362527499258
$arr[9,5] = 3.601663961009424
$arr[9,6] = 0
$arr[9,7] = 17.80965308857159
$arr[9,8] = 7.503458201869291
$arr[9,9] = 0
$arr[9,10] = 0
$arr[9,11] = 19.79762216296736
$arr[9,12] = 16.68065842654093
$arr[9,13] = 19.16913916008912
$arr[10,0] = 13.39404025393302
$arr[10,1] = 9.471000807609055
$arr[10,2] = 0
$arr[10,3] = 22.02905412376154
$arr[10,4] = 38.48873211151496
$arr[10,5] = 3.601192718652287
$arr[10,6] = 0
$arr[10,7] = 17.78819034344317
$arr[10,8] = 7.499394389138343
$arr[10,9] = 0
$arr[10,10] = 0
$arr[10,11] = 19.85895998516512
$arr[10,12] = 16.66965034288686
$arr[10,13] = 19.1652337616231
$arr[11,0] = 13.36522586390503
$arr[11,1] = 9.449628947453828
$arr[11,2] = 0
$arr[11,3] = 22.02515143889718
$arr[11,4] = 38.47895347603666
$arr[11,5] = 3.601293808621977
$arr[11,6] = 0
$arr[11,7] = 17.7927836566317
$arr[11,8] = 7.500266380645157
$arr[11,9] = 0
$arr[11,10] = 0
$arr[11,11] = 19.84574696531326
$arr[11,12] = 16.67200955632614
$arr[11,13] = 19.16603803116208
$arr[12,0] = 13.27082379888946
$arr[12,1] = 9.379557436519907
$arr[12,2] = 0
$arr[12,3] = 22.01252248136883
$arr[12,4] = 38.44730732995611
$arr[12,5] = 3.60162501107665
$arr[12,6] = 0
$arr[12,7] = 17.80787423862213
$arr[12,8] = 7.503122416586636
$arr[12,9] = 0
$arr[12,10] = 0
$arr[12,11] = 19.80266712054292
$arr[12,12] = 16.67974756837721
$arr[12,13] = 19.16880126865879
$arr[13,0] = 13.21266888543048
$arr[13,1] = 9.336350041264744
$arr[13,2] = 0
$arr[13,3] = 22.00486276870725
$arr[13,4] = 38.42811122577463
$arr[13,5] = 3.601829055630648
$arr[13,6] = 0
$arr[13,7] = 17.81720273968651
$arr[13,8] = 7.504881264962537
$arr[13,9] = 0
$arr[13,10] = 0
$arr[13,11] = 19.77628855765789
$arr[13,12] = 16.68452122544777
$arr[13,13] = 19.17060162750984
$arr[14,0] = 12.8743818101308
$arr[14,1] = 9.084363639774946
$arr[14,2] = 0
$arr[14,3] = 21.96213784694235
$arr[14,4] = 38.32100877622017
$arr[14,5] = 3.603016372334778
$arr[14,6] = 0
$arr[14,7] = 17.87196271452571
$arr[14,8] = 7.515105515133532
$arr[14,9] = 0
$arr[14,10] = 0
$arr[14,11] = 19.62529920159169
$arr[14,12] = 16.7123971857574
$arr[14,13] = 19.18256907676437
$arr[15,0] = 12.66242495683602
$arr[15,1] = 8.92587991651634
$arr[15,2] = 0
$arr[15,3] = 21.93697832782359
$arr[15,4] = 38.25791602942537
$arr[15,5] = 3.603760868622515
$arr[15,6] = 0
$arr[15,7] = 17.90671592517457
$arr[15,8] = 7.521507423296958
$arr[15,9] = 0
$arr[15,10] = 0
$arr[15,11] = 19.53286832055679
$arr[15,12] = 16.72996273750997
$arr[15,13] = 19.19138111059254
$arr[16,0] = 12.53890880106909
$arr[16,1] = 8.833297046289761
$arr[16,2] = 0
$arr[16,3] = 21.92289610909799
$arr[16,4] = 38.22259437395717
$arr[16,5] = 3.60419502026158
$arr[16,6] = 0
$arr[16,7] = 17.92713040679523
$arr[16,8] = 7.525237390994826
$arr[16,9] = 0
$arr[16,10] = 0
$arr[16,11] = 19.47979236645537
$arr[16,12] = 16.7402368977397
$arr[16,13] = 19.19698897965667
$arr[17,0] = 12.49681519945706
$arr[17,1] = 8.801705511944164
$arr[17,2] = 0
$arr[17,3] = 21.91819514807989
$arr[17,4] = 38.21080199260954
$arr[17,5] = 3.604343037624481
$arr[17,6] = 0
$arr[17,7] = 17.93411537493933
$arr[17,8] = 7.526508507083397
$arr[17,9] = 0
$arr[17,10] = 0
$arr[17,11] = 19.4618383445831
$arr[17,12] = 16.74374492645681
$arr[17,13] = 19.19898025840395
$arr[18,0] = 12.68515470473596
$arr[18,1] = 8.942898621766622
$arr[18,2] = 0
$arr[18,3] = 21.93961641415405
$arr[18,4] = 38.2645323908708
$arr[18,5] = 3.603681001585448
$arr[18,6] = 0
$arr[18,7] = 17.90297234504556
$arr[18,8] = 7.520820988784095
$arr[18,9] = 0
$arr[18,10] = 0
$arr[18,11] = 19.54269900975035
$arr[18,12] = 16.7280751695096
$arr[18,13] = 19.19038720425147
$arr[19,0] = 13.29862104669274
$arr[19,1] = 9.400198881679644
$arr[19,2] = 0
$arr[19,3] = 22.01621611336993
$arr[19,4] = 38.45656342449744
$arr[19,5] = 3.601527484413658
$arr[19,6] = 0
$arr[19,7] = 17.80342403268001
$arr[19,8] = 7.502281562101842
$arr[19,9] = 0
$arr[19,10] = 0
$arr[19,11] = 19.81531890829952
$arr[19,12] = 16.67746766273539
$arr[19,13] = 19.16796716986498
$arr[20,0] = 13.6848054227985
$arr[20,1] = 9.686258191173232
$arr[20,2] = 0
$arr[20,3] = 22.06968242365237
$arr[20,4] = 38.59050670705757
$arr[20,5] = 3.60017259793329
$arr[20,6] = 0
$arr[20,7] = 17.74216931837407
$arr[20,8] = 7.490587954599446
$arr[20,9] = 0
$arr[20,10] = 0
$arr[20,11] = 19.99394074572108
$arr[20,12] = 16.64591059950563
$arr[20,13] = 19.15813664277569
$arr[21,0] = 13.48006151722307
$arr[21,1] = 9.534759531293723
$arr[21,2] = 0
$arr[21,3] = 22.04083796937599
$arr[21,4] = 38.51825543246959
$arr[21,5] = 3.600890931778052
$arr[21,6] = 0
$arr[21,7] = 17.7745129851265
$arr[21,8] = 7.496790461783003
$arr[21,9] = 0
$arr[21,10] = 0
$arr[21,11] = 19.89858185350632
$arr[21,12] = 16.66261450962515
$arr[21,13] = 19.16294133749094
$arr[22,0] = 12.67488374867625
$arr[22,1] = 8.935209037866272
$arr[22,2] = 0
$arr[22,3] = 21.93842254376051
$arr[22,4] = 38.26153816897973
$arr[22,5] = 3.603717090372347
$arr[22,6] = 0
$arr[22,7] = 17.90466346466977
$arr[22,8] = 7.521131171841108
$arr[22,9] = 0
$arr[22,10] = 0
$arr[22,11] = 19.5382543527066
$arr[22,12] = 16.72892799229871
$arr[22,13] = 19.19083486209832
$arr[23,0] = 11.74544579990607
$arr[23,1] = 8.234038343645697
$arr[23,2] = 0
$arr[23,3] = 21.84284551938443
$arr[23,4] = 38.0216996691376
$arr[23,5] = 3.606992738966994
$arr[23,6] = 0
$arr[23,7] = 18.06130403838104
$arr[23,8] = 7.549213969975484
$arr[23,9] = 0
$arr[23,10] = 0
$arr[23,11] = 19.15314051141221
$arr[23,12] = 16.80699523249699
$arr[23,13] = 19.24152923874994

$ws.Range("B2:O25").Value = $arr

